$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Define the three new character styles referenced by the edited runs.
# ---------------------------------------------------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2. Apply "GaNStyle" to every run containing the campaign dates sentence.
#    (occurs 4 times throughout the document)
# ---------------------------------------------------------------------------

$datesText = "2022: Datumi kampanje za opazovanje Herkulovo ozvezdje: 13.-22. junij, 12.-21. julij, 10.-19. avgust"
$rng = $d.Content
while ($rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# ---------------------------------------------------------------------------
# 3. Apply "GaNParagraph" to the run with the "Sodelujete v svetovni..." text.
# ---------------------------------------------------------------------------

$paragraphText = "Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega Herkulovo ozvezdje na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom."
$rng2 = $d.Content
if ($rng2.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
}

# ---------------------------------------------------------------------------
# 4. Apply "GaNLinks" to the run with the "Jenik Hollan, CzechGlobe..." text.
# ---------------------------------------------------------------------------

$linksText = "Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
if ($rng3.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
}

Write-Host "Styles created and applied."
